$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)
Write-Output $sh.Name
$sh.TextFrame.TextRange.Text = "10/13/2020"
